$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.061.13"
$ws.Range("E2").Value = "  -3.25%  "
$ws.Range("D3").Value = "1.644.94"
$ws.Range("E3").Value = "  -5.03%  "
$ws.Range("D4").Value = "'0.9984"
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").Value = "'233.32"
$ws.Range("E5").Value = "  -5.03%  "
$ws.Range("D6").Value = "'0.9994"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").Value = "'0.4735"
$ws.Range("E7").Value = "  -5.43%  "
$ws.Range("D8").Value = "'0.2562"
$ws.Range("E8").Value = "  -5.63%  "
$ws.Range("D9").Value = "'0.06084"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D10").Value = "'0.07031"
$ws.Range("E10").Value = "  -2.90%  "
$ws.Range("D11").Value = "1.640.90"
$ws.Range("E11").Value = "  -5.38%  "
$ws.Range("D12").Value = "'14.32"
$ws.Range("E12").Value = "  -5.32%  "
$ws.Range("D13").Value = "'4.301"
$ws.Range("E13").Value = "  -9.15%  "
$ws.Range("D14").Value = "'0.5739"
$ws.Range("E14").Value = "  -12.29%  "
$ws.Range("D15").Value = "'73.40"
$ws.Range("E15").Value = "  -4.71%  "
$ws.Range("D16").Value = "'0.9990"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "'0.9995"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "25.061.40"
$ws.Range("E18").Value = "  -3.41%  "
$ws.Range("D19").Value = "'11.27"
$ws.Range("E19").Value = "  -4.56%  "
$ws.Range("D20").Value = "'0.000006581"
$ws.Range("E20").Value = "  -3.05%  "
$ws.Range("D21").Value = "1.849.85"
$ws.Range("E21").Value = "  -5.72%  "
$ws.Range("D22").Value = "'4.294"
$ws.Range("E22").Value = "  -6.65%  "
$ws.Range("D23").Value = "'8.456"
$ws.Range("E23").Value = "  -3.48%  "
$ws.Range("D24").Value = "'5.235"
$ws.Range("E24").Value = "  -3.05%  "
$ws.Range("D25").Value = "'133.96"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").Value = "'14.90"
$ws.Range("E26").Value = "  -2.23%  "
$ws.Range("D27").Value = "'1.374"
$ws.Range("E27").Value = "  -3.08%  "
$ws.Range("D28").Value = "'103.39"
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("D29").Value = "'1.635"
$ws.Range("E29").Value = "  -8.13%  "
$ws.Range("D30").Value = "'3.884"
$ws.Range("E30").Value = "  -2.01%  "
$ws.Range("D31").Value = "'0.07573"
$ws.Range("E31").Value = "  -6.35%  "
$ws.Range("D32").Value = "'3.519"
$ws.Range("E32").Value = "  -4.79%  "
$ws.Range("D33").Value = "'0.9986"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("D34").Value = "'0.04236"
$ws.Range("E34").Value = "  -10.50%  "
$ws.Range("D35").Value = "'2.575"
$ws.Range("E35").Value = "  -3.38%  "
$ws.Range("D36").Value = "'0.9318"
$ws.Range("E36").Value = "  -6.16%  "
$ws.Range("D37").Value = "'0.5897"
$ws.Range("E37").Value = "  -2.62%  "
$ws.Range("D38").Value = "'2.585"
$ws.Range("E38").Value = "  -5.88%  "
$ws.Range("D39").Value = "'0.8566"
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("D40").Value = "'0.9990"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").Value = "'0.01482"
$ws.Range("E41").Value = "  -7.63%  "
$ws.Range("D42").Value = "'98.12"
$ws.Range("E42").Value = "  -2.15%  "
$ws.Range("D43").Value = "'1.776"
$ws.Range("E43").Value = "  -8.10%  "
$ws.Range("D44").Value = "'0.3671"
$ws.Range("E44").Value = "  -6.19%  "
$ws.Range("D45").Value = "'4.619"
$ws.Range("E45").Value = "  -7.56%  "
$ws.Range("D46").Value = "'0.1094"
$ws.Range("E46").Value = "  -6.92%  "
$ws.Range("D47").Value = "'0.05204"
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D48").Value = "'6.062"
$ws.Range("E48").Value = "  -4.21%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'28.79"
$ws.Range("E49").Value = "  -6.50%  "
$ws.Range("B50").Value = "TrueUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D50").Value = "'0.9990"
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("D51").Value = "'0.9993"
$ws.Range("E51").Value = "  -0.36%  "
